$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.646.27"
$ws.Range("E2").Value = "  +0.37%  "

$ws.Range("D3").Value = "1.596.95"
$ws.Range("E3").Value = "  +0.84%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.75"
$ws.Range("E5").Value = "  +0.03%  "

$ws.Range("E6").Value = "  +1.27%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("E8").Value = "  -0.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.246"
$ws.Range("E9").Value = "  -1.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.58"
$ws.Range("E10").Value = "  +0.51%  "

$ws.Range("E11").Value = "  +1.20%  "

$ws.Range("D12").Value = "1.825.15"
$ws.Range("E12").Value = "  +1.11%  "

$ws.Range("D13").Value = "1.625.28"
$ws.Range("E13").Value = "  +2.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.01"
$ws.Range("E14").Value = "  -0.67%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.519"
$ws.Range("E15").Value = "  -1.55%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.72"
$ws.Range("E16").Value = "  +0.86%  "

$ws.Range("D17").Value = "26.645.17"
$ws.Range("E17").Value = "  +0.28%  "

$ws.Range("D18").Value = "0.0₃0727"
$ws.Range("E18").Value = "  -0.22%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "208.54"
$ws.Range("E19").Value = "  +0.07%  "

$ws.Range("E20").Value = "  +0.04%  "

$ws.Range("E21").Value = "  +0.72%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.24"
$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.31"
$ws.Range("E23").Value = "  -3.40%  "

$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.47"
$ws.Range("E25").Value = "  -0.65%  "

$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.21"
$ws.Range("E27").Value = "  -2.91%  "

$ws.Range("E28").Value = "  +2.07%  "

$ws.Range("E29").Value = "  -0.22%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0505"
$ws.Range("E30").Value = "  +0.66%  "

$ws.Range("E31").Value = "  +0.15%  "

$ws.Range("E32").Value = "  -1.08%  "

$ws.Range("E33").Value = "  -3.10%  "

$ws.Range("E34").Value = "  -0.23%  "

$ws.Range("D35").Value = "1.290.50"
$ws.Range("E35").Value = "  -2.01%  "

$ws.Range("E36").Value = "  +0.56%  "

$ws.Range("E37").Value = "  -1.10%  "

$ws.Range("E38").Value = "  -1.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.844"
$ws.Range("E39").Value = "  +2.90%  "

$ws.Range("E40").Value = "  +0.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.40"
$ws.Range("E41").Value = "  +2.02%  "

$ws.Range("E42").Value = "  +1.39%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.786"
$ws.Range("E43").Value = "  +0.27%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.67"
$ws.Range("E44").Value = "  +1.15%  "

$ws.Range("D45").Value = "1.736.77"
$ws.Range("E45").Value = "  +1.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.899"
$ws.Range("E46").Value = "  +8.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.97"
$ws.Range("E47").Value = "  +1.20%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.60"
$ws.Range("E48").Value = "  -0.98%  "

$ws.Range("E49").Value = "  +1.79%  "

$ws.Range("E50").Value = "  -0.28%  "

$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₇0967"
$ws.Range("E51").Value = "  -0.68%  "
